# Update the cryptocurrency price/volume table with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=rank, B=Coin, C=Link, D=Price, E=Volume(1h)

# Helper: write a text value to a cell without letting Excel auto-convert
# digit-only strings (e.g. "1.00", "596.61") into numbers, and without
# leaving a residual style on the cell once written.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 43 and 44 swap positions (OKB <-> ONDO) plus updated figures.
$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D43") "1.24"
$ws.Range("E43").Value = "  +3.97%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D44") "42.62"
$ws.Range("E44").Value = "  +0.56%  "

# Price (D) and Volume(1h) (E) updates for the remaining rows.
$ws.Range("D2").Value = "64.846.78"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "3.521.26"
$ws.Range("E3").Value = "  +0.01%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue $ws.Range("D5") "596.61"
$ws.Range("E5").Value = "  +0.50%  "

Set-TextValue $ws.Range("D6") "134.18"
$ws.Range("E6").Value = "  -1.63%  "

$ws.Range("D7").Value = "3.521.35"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.62%  "

$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("E11").Value = "  +3.92%  "

$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "4.119.38"
$ws.Range("E13").Value = "  +0.04%  "

Set-TextValue $ws.Range("D14") "27.34"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "3.517.82"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "64.921.14"
$ws.Range("E18").Value = "  +0.01%  "

Set-TextValue $ws.Range("D19") "9.90"
$ws.Range("E19").Value = "  -1.07%  "

Set-TextValue $ws.Range("D20") "14.44"
$ws.Range("E20").Value = "  +2.06%  "

Set-TextValue $ws.Range("D21") "5.70"
$ws.Range("E21").Value = "  -1.99%  "

Set-TextValue $ws.Range("D22") "391.24"
$ws.Range("E22").Value = "  +0.82%  "

Set-TextValue $ws.Range("D23") "0.576"
$ws.Range("E23").Value = "  +0.86%  "

$ws.Range("D24").Value = "3.659.99"
$ws.Range("E24").Value = "  -0.01%  "

Set-TextValue $ws.Range("D25") "74.10"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("E28").Value = "  +18.30%  "

Set-TextValue $ws.Range("D29") "7.79"
$ws.Range("E29").Value = "  +1.97%  "

Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("E31").Value = "  +1.82%  "

Set-TextValue $ws.Range("D32") "8.40"
$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("D33").Value = "3.521.94"
$ws.Range("E33").Value = "  -0.37%  "

Set-TextValue $ws.Range("D34") "24.08"
$ws.Range("E34").Value = "  +1.60%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +1.41%  "

Set-TextValue $ws.Range("D37") "5.23"
$ws.Range("E37").Value = "  +6.07%  "

Set-TextValue $ws.Range("D38") "1.58"
$ws.Range("E38").Value = "  +2.18%  "

Set-TextValue $ws.Range("D39") "169.01"
$ws.Range("E39").Value = "  -0.18%  "

Set-TextValue $ws.Range("D40") "6.85"
$ws.Range("E40").Value = "  +0.66%  "

Set-TextValue $ws.Range("D41") "0.0821"
$ws.Range("E41").Value = "  +3.11%  "

Set-TextValue $ws.Range("D42") "0.821"
$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("E45").Value = "  +0.06%  "

Set-TextValue $ws.Range("D46") "25.25"
$ws.Range("E46").Value = "  -4.62%  "

Set-TextValue $ws.Range("D47") "4.43"
$ws.Range("E47").Value = "  +0.53%  "

$ws.Range("E48").Value = "  -0.44%  "

$ws.Range("E49").Value = "  +1.44%  "

$ws.Range("D50").Value = "2.391.41"
$ws.Range("E50").Value = "  -0.26%  "

Set-TextValue $ws.Range("D51") "0.900"
$ws.Range("E51").Value = "  +7.23%  "
